# CSP Info.xlsx - "Add files via upload" edit
#
# Row 11 (Logan Fabris) is missing values for the "Favorite Ice Cream"
# and "Favorite Pizza Toppings" table columns (E11 / F11). Fill them in,
# which also appends two new shared-string entries ("Chocolate",
# "Sausage") to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "Chocolate"
$ws.Range("F11").Value = "Sausage"

# Reflect the cursor/selection position at save time.
[void]$ws.Range("H11").Select()
